$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.308.35'
$ws.Range("E2").Value = '  -6.98%  '
$ws.Range("D3").Value = '3.736.97'
$ws.Range("E3").Value = '  -5.95%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '581.04'
$ws.Range("E5").Value = '  -5.15%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.16'
$ws.Range("E6").Value = '  +3.62%  '
$ws.Range("D7").Value = '3.726.23'
$ws.Range("E7").Value = '  -5.96%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.634'
$ws.Range("E8").Value = '  -6.90%  '
$ws.Range("E9").Value = '  +0.22%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.718'
$ws.Range("E10").Value = '  -6.08%  '
$ws.Range("E11").Value = '  -9.79%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '53.81'
$ws.Range("E12").Value = '  -3.86%  '
$ws.Range("E13").Value = '  -10.34%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.80'
$ws.Range("E14").Value = '  -3.66%  '
$ws.Range("D15").Value = '4.337.33'
$ws.Range("E15").Value = '  -5.93%  '
$ws.Range("D16").Value = '3.738.36'
$ws.Range("E16").Value = '  -6.03%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.51'
$ws.Range("E17").Value = '  -4.37%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.16'
$ws.Range("E18").Value = '  -7.05%  '
$ws.Range("E19").Value = '  -7.02%  '
$ws.Range("E20").Value = '  -2.79%  '
$ws.Range("D21").Value = '68.138.64'
$ws.Range("E21").Value = '  -6.95%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '412.67'
$ws.Range("E22").Value = '  -6.05%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.57'
$ws.Range("E23").Value = '  -5.98%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '89.07'
$ws.Range("E24").Value = '  -6.90%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.10'
$ws.Range("E25").Value = '  -7.96%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.96'
$ws.Range("E26").Value = '  -8.58%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.85'
$ws.Range("E27").Value = '  -1.50%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.85'
$ws.Range("E28").Value = '  -5.42%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.00'
$ws.Range("E29").Value = '  +0.73%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.63'
$ws.Range("E30").Value = '  -8.12%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.00'
$ws.Range("E31").Value = '  +2.63%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '33.19'
$ws.Range("E32").Value = '  -7.98%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '12.79'
$ws.Range("E33").Value = '  -7.69%  '
$ws.Range("E34").Value = '  -8.13%  '
$ws.Range("B35").Value = 'Bittensor'
$ws.Range("C35").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '616.95'
$ws.Range("E35").Value = '  -4.65%  '
$ws.Range("B36").Value = 'OKB'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '65.61'
$ws.Range("E36").Value = '  -6.72%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '43.88'
$ws.Range("E37").Value = '  -8.30%  '
$ws.Range("D38").Value = '0.0₃0927'
$ws.Range("E38").Value = '  -11.93%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.402'
$ws.Range("E39").Value = '  -6.24%  '
$ws.Range("E40").Value = '  +0.15%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.04%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.20'
$ws.Range("E42").Value = '  -0.88%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.137'
$ws.Range("E43").Value = '  -5.72%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.09'
$ws.Range("E44").Value = '  -8.69%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0448'
$ws.Range("E45").Value = '  -7.62%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.64'
$ws.Range("E46").Value = '  +2.87%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.45'
$ws.Range("E47").Value = '  -10.58%  '
$ws.Range("E48").Value = '  -8.19%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.72'
$ws.Range("E49").Value = '  -15.50%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.19'
$ws.Range("E50").Value = '  -6.82%  '
$ws.Range("D51").Value = '2.733.74'
$ws.Range("E51").Value = '  -2.34%  '
